$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B-E and G, rows 2-9 (regen sval data to filter save games)
$data = @{
    2 = @{ B = 1.505614041169197;  C = 1.65323645889881;  D = 0.1529057820181812; E = 0.4998867070740569; G = 3.811642989160245 }
    3 = @{ B = 1.505614041169197;  C = 1.65323645889881;  D = 0.7127328510149897; E = 0.4998867070740569; G = 4.371470058157054 }
    4 = @{ B = 3.182878228561681;  C = 1.65323645889881;  D = 3.082599426703578;  E = 0.4998867070740569; G = 8.418600821238126 }
    5 = @{ B = 1.505614041169197;  C = 1.65323645889881;  D = 0.1529057820181812; E = 0.4998867070740569; G = 3.811642989160245 }
    6 = @{ B = 3.182878228561681;  C = 1.65323645889881;  D = 3.082599426703578;  E = 6.48142807727062;   G = 14.40014219143469 }
    7 = @{ B = 1.505614041169197;  C = 1.65323645889881;  D = 0.7127328510149897; E = 0.4998867070740569; G = 4.371470058157054 }
    8 = @{ B = 0.1554434735375247; C = 1.65323645889881;  D = 0.7127328510149897; E = 246.9852506941017;  G = 249.506663477553 }
    9 = @{ B = 1.505614041169197;  C = 0.3375848360084654; D = 0.7127328510149897; E = 0.4998867070740569; G = 3.055818435266709 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
